$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 4.916666666666667
$ws.Range("C2").Value = 5

$ws.Range("B3").Value = 4.916666666666667
$ws.Range("C3").Value = 5

$ws.Range("B4").Value = 3.958333333333333
$ws.Range("C4").Value = 4

$ws.Range("B5").Value = 8.041666666666666
$ws.Range("C5").Value = 9

$ws.Range("B6").Value = 18.16666666666667
$ws.Range("C6").Value = 19
